$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new range to be formatted as Text so that date-like and time-like
# strings (e.g. "2024-06-10", "11:25:14", "0:00:01") are stored verbatim as
# text rather than being auto-converted into date/time serial numbers.
$newRange = $ws.Range("A171:H175")
$newRange.NumberFormat = "@"

# Row 171
$ws.Range("A171").Value = "WV50 FILTER"
$ws.Range("B171").Value = "NOK Soldadura metal"
$ws.Range("C171").Value = "2024-06-10"
$ws.Range("D171").Value = "11:25:14"
$ws.Range("E171").Value = "Mañana"
$ws.Range("F171").Value = "11:25:15"
$ws.Range("G171").Value = "0:00:01"
$ws.Range("H171").Value = "-0.00 minutos"

# Row 172
$ws.Range("A172").Value = "WV50 FILTER"
$ws.Range("B172").Value = "Traza"
$ws.Range("C172").Value = "2024-06-10"
$ws.Range("D172").Value = "11:25:20"
$ws.Range("E172").Value = "Mañana"
$ws.Range("F172").Value = "11:25:21"
$ws.Range("G172").Value = "0:00:01"
$ws.Range("H172").Value = "0.02 minutos"

# Row 173
$ws.Range("A173").Value = "WV50 FILTER"
$ws.Range("B173").Value = "Fallo cámara ferrite"
$ws.Range("C173").Value = "2024-06-10"
$ws.Range("D173").Value = "11:25:53"
$ws.Range("E173").Value = "Mañana"
$ws.Range("F173").Value = "11:25:57"
$ws.Range("G173").Value = "0:00:04"
$ws.Range("H173").Value = "0.11 minutos"

# Row 174
$ws.Range("A174").Value = "WC47 NACP"
$ws.Range("B174").Value = "No coge placa"
$ws.Range("C174").Value = "2024-06-10"
$ws.Range("D174").Value = "11:32:58"
$ws.Range("E174").Value = "Mañana"
$ws.Range("F174").Value = "11:32:59"
$ws.Range("G174").Value = "0:00:01"
$ws.Range("H174").Value = "-0.00 minutos"

# Row 175
$ws.Range("A175").Value = "WC47 NACP"
$ws.Range("B175").Value = "No pone tornillo"
$ws.Range("C175").Value = "2024-06-10"
$ws.Range("D175").Value = "11:36:34"
$ws.Range("E175").Value = "Mañana"
$ws.Range("F175").Value = "11:36:36"
$ws.Range("G175").Value = "0:00:02"
$ws.Range("H175").Value = "-0.01 minutos"

